# Friday (2/26) time tracking.
# Week-2 row (row 4, week starting 2021-02-22): record minutes worked on
# Tuesday (column D) and Friday (column G). The shared "Total Minutes"
# (J) and "Hours" (K) formulas, plus the J49/K49 grand totals, recalc
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 15
$ws.Range("G4").Value = 55

# Leave the cursor parked on G5, matching the saved selection state.
$ws.Range("G5").Select()
